# ENR for noise figure meter automation / qLoadENR table additions.
# Apply to the "BigBookOfVariableNames" worksheet:
#   1. Insert a blank row above row 2 (shifts existing rows 2-21 down to 3-22).
#   2. Append new "calibration" / "results" / "Model bands" rows at 26-29.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BigBookOfVariableNames")

# 1. Insert a new blank row at row 2, pushing everything below down by one.
$ws.Rows.Item(2).Insert()

# 2. Add the new rows describing calibration / results / model bands.
#    Values are written in the same order the original author typed them
#    (column A top-to-bottom first, then column B) so new shared-string
#    entries land in the same order as the source workbook.
$ws.Range("A26").Value = "calibration"
$ws.Range("A28").Value = "results"
$ws.Range("A29").Value = "Model bands"

$ws.Range("B28").Value = "Input Frequency, Band, temperature, timestamp, tech name, bench, calibration files used (timestamp)"
$ws.Range("B26").Value = "frequency(input with corresponding output by step) , input network correction, output network correction, ""happy range"", timestamp"
$ws.Range("B29").Value = "frequencies, powers, bands, LO's, step sizes per band,"
$ws.Range("B27").Value = "calculate output frequency save and sort calibration table before save. Write a cal table plot function."

# 3. Restore the selection/active cell to B27 as in the saved file.
$ws.Range("B27").Select()
